$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 = state_id: Type smallint -> varchar, Sample "6" -> "06" (kept as text)
$ws.Range("C4").Value = "varchar"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "06"

# Row 5 = state_name: Sample stays "California" (no value change required)
$ws.Range("E5").Value = "California"

# Update selection to E3 to match the saved view state
$ws.Range("E3").Select()

# Mark the "number stored as text" warning on E4 as ignored
$ws.Range("E4").Errors.Item([Microsoft.Office.Interop.Excel.XlErrorChecks]::xlNumberAsText).Ignore = $true
